# B6-PowerPoint.pptx edit
#
# 1) Three table graphic-frames (on what were originally slides 14, 15 and
#    16) get their table style switched from
#    {3EEAD695-FF23-497C-81A7-C03473788C5A} to
#    {DD6D9A8A-93EF-4155-BC3A-A262D7E38635}.
#
# 2) The deck's theme is swapped: the palette that is actually applied to
#    the slide master (currently the "Integral" / "Red Violet" scheme)
#    becomes the stock "Office Theme" palette (the one that, before the
#    edit, only the notes master pointed at).  The two themes share an
#    identical font scheme / format scheme, so the only thing that
#    actually differs between them is the 12-colour scheme - that's what
#    we reapply here.

$p = $ppt.ActivePresentation

# --- 1) Retarget the three affected tables to the new table style ----
$targetStyle = "{DD6D9A8A-93EF-4155-BC3A-A262D7E38635}"
foreach ($idx in 14, 15, 16) {
    $slide = $p.Slides.Item($idx)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyle)
        }
    }
}

# --- 2) Swap the active theme's colour scheme to the Office palette --
$colors = $p.SlideMaster.Theme.ThemeColorScheme
$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72

Write-Host "table styles + theme colours updated"
